$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sending cluster changes from FAPs to ECs; other values recomputed
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ngf"
$ws.Range("C2").Value = "Sorcs3"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05916633333333333
$ws.Range("H2").Value = 0.177499
$ws.Range("I2").Value = 0.005840587905374044
$ws.Range("J2").Value = 0.005840587905374044
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.08256033333333333
$ws.Range("N2").Value = 0.247681
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.004884792202111111
$ws.Range("R2").Value = 0.043963129819
$ws.Range("S2").Value = 0.005840587905374044
$ws.Range("T2").Value = 0.005840587905374044

# Row 3: sending cluster changes from sCs to FAPs; other values recomputed
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ngf"
$ws.Range("C3").Value = "Sorcs3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.748364666666667
$ws.Range("H3").Value = 17.245094
$ws.Range("I3").Value = 0.567448196572592
$ws.Range("J3").Value = 0.567448196572592
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.08256033333333333
$ws.Range("N3").Value = 0.247681
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.4745869030015556
$ws.Range("R3").Value = 4.271282127014
$ws.Range("S3").Value = 0.567448196572592
$ws.Range("T3").Value = 0.567448196572592

# Row 4 (new): sending cluster sCs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ngf"
$ws.Range("C4").Value = "Sorcs3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.322670666666667
$ws.Range("H4").Value = 12.968012
$ws.Range("I4").Value = 0.4267112155220338
$ws.Range("J4").Value = 0.4267112155220338
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.08256033333333333
$ws.Range("N4").Value = 0.247681
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.3568811311302222
$ws.Range("R4").Value = 3.211930180172
$ws.Range("S4").Value = 0.4267112155220338
$ws.Range("T4").Value = 0.4267112155220338
